# Generate Report for handoff
#
# This script updates the three worksheets (Overview, zh-cn, de-de) of the
# localization-status workbook so that:
#   - the "865ddbfb-...md" row moves from position 2 to position 4 and its
#     status becomes "Ready for handoff" (with refreshed handoff/handback
#     datetimes),
#   - the "ffff733116a9-...md" row moves from position 3 to position 2,
#   - the "ffffff5c4452c2-...md" row moves from position 4 to position 3.
#
# Hyperlink target addresses are intentionally left untouched (only the
# cell text / hyperlink display text is updated), matching the underlying
# OOXML relationships which are not modified by this change.

$wb = $excel.ActiveWorkbook

function Set-CellAndHyperlink {
    param(
        $ws,
        [string]$addr,
        [string]$newValue
    )
    $ws.Range($addr).Value2 = $newValue
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq ('$' + ($addr -replace '(\d+)', '$$$1'))) {
            $hl.TextToDisplay = $newValue
        }
    }
}

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-CellAndHyperlink $wsOverview "A2" "ffff733116a9-d2b7-4004-8f23-79f05e03e5b8.md"
Set-CellAndHyperlink $wsOverview "A3" "ffffff5c4452c2-b997-4aa4-a951-c8bed0ec560f.md"
Set-CellAndHyperlink $wsOverview "A4" "865ddbfb-e348-460f-859b-030bdc325eb9.md"

$wsOverview.Range("B4").Value2 = "Ready for handoff"
$wsOverview.Range("C4").Value2 = "Ready for handoff"

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-CellAndHyperlink $wsZhCn "A2" "ffff733116a9-d2b7-4004-8f23-79f05e03e5b8.md"
Set-CellAndHyperlink $wsZhCn "C2" "4f872ba6-5f5e-4fed-ae4b-08952f4241ec.3582d9695a9dd21afa17630a94680e3c8b126d09.zh-cn.xlf"
$wsZhCn.Range("D2").Value2 = "2016-01-25 04:01:35"
Set-CellAndHyperlink $wsZhCn "E2" "4f872ba6-5f5e-4fed-ae4b-08952f4241ec.md"
Set-CellAndHyperlink $wsZhCn "F2" "4f872ba6-5f5e-4fed-ae4b-08952f4241ec.3582d9695a9dd21afa17630a94680e3c8b126d09.zh-cn.xlf"
$wsZhCn.Range("G2").Value2 = "2016-01-25 04:02:19"

Set-CellAndHyperlink $wsZhCn "A3" "ffffff5c4452c2-b997-4aa4-a951-c8bed0ec560f.md"

Set-CellAndHyperlink $wsZhCn "A4" "865ddbfb-e348-460f-859b-030bdc325eb9.md"
$wsZhCn.Range("B4").Value2 = "Ready for handoff"
Set-CellAndHyperlink $wsZhCn "C4" "865ddbfb-e348-460f-859b-030bdc325eb9.d521806228e9e9c94733ee6eb4c8fc973b7a4bff.zh-cn.xlf"
$wsZhCn.Range("D4").Value2 = "2016-01-25 04:05:37"
Set-CellAndHyperlink $wsZhCn "E4" "865ddbfb-e348-460f-859b-030bdc325eb9.md"
Set-CellAndHyperlink $wsZhCn "F4" "865ddbfb-e348-460f-859b-030bdc325eb9.d521806228e9e9c94733ee6eb4c8fc973b7a4bff.zh-cn.xlf"
$wsZhCn.Range("G4").Value2 = "2016-01-25 04:04:34"

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-CellAndHyperlink $wsDeDe "A2" "ffff733116a9-d2b7-4004-8f23-79f05e03e5b8.md"
Set-CellAndHyperlink $wsDeDe "C2" "4f872ba6-5f5e-4fed-ae4b-08952f4241ec.3582d9695a9dd21afa17630a94680e3c8b126d09.de-de.xlf"
$wsDeDe.Range("D2").Value2 = "2016-01-25 04:01:45"
Set-CellAndHyperlink $wsDeDe "E2" "4f872ba6-5f5e-4fed-ae4b-08952f4241ec.md"
Set-CellAndHyperlink $wsDeDe "F2" "4f872ba6-5f5e-4fed-ae4b-08952f4241ec.3582d9695a9dd21afa17630a94680e3c8b126d09.de-de.xlf"
$wsDeDe.Range("G2").Value2 = "2016-01-25 04:02:36"

Set-CellAndHyperlink $wsDeDe "A3" "ffffff5c4452c2-b997-4aa4-a951-c8bed0ec560f.md"

Set-CellAndHyperlink $wsDeDe "A4" "865ddbfb-e348-460f-859b-030bdc325eb9.md"
$wsDeDe.Range("B4").Value2 = "Ready for handoff"
Set-CellAndHyperlink $wsDeDe "C4" "865ddbfb-e348-460f-859b-030bdc325eb9.d521806228e9e9c94733ee6eb4c8fc973b7a4bff.de-de.xlf"
$wsDeDe.Range("D4").Value2 = "2016-01-25 04:05:47"
Set-CellAndHyperlink $wsDeDe "E4" "865ddbfb-e348-460f-859b-030bdc325eb9.md"
Set-CellAndHyperlink $wsDeDe "F4" "865ddbfb-e348-460f-859b-030bdc325eb9.d521806228e9e9c94733ee6eb4c8fc973b7a4bff.de-de.xlf"
$wsDeDe.Range("G4").Value2 = "2016-01-25 04:04:52"
